$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need to be forced to Text
# format first, otherwise Excel auto-converts them to a Number and the
# original formatted text (e.g. trailing zeros) is lost.

$ws.Range('D2').Value = '27.662.36'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.581.63'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.91'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.20'
$ws.Range('E8').Value = '  -5.32%  '
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0590'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').Value = '1.807.88'
$ws.Range('E12').Value = '  -3.28%  '
$ws.Range('D13').Value = '1.581.23'
$ws.Range('E13').Value = '  -3.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.84'
$ws.Range('E14').Value = '  -4.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  -6.60%  '
$ws.Range('D16').Value = '27.626.80'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.13'
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.12'
$ws.Range('E18').Value = '  -4.88%  '
$ws.Range('D19').Value = '0.0₃0690'
$ws.Range('E19').Value = '  -4.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.29'
$ws.Range('E20').Value = '  -7.38%  '
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('E22').Value = '  -5.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.47'
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.45'
$ws.Range('E25').Value = '  -1.94%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.74'
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.08'
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('E29').Value = '  -4.24%  '
$ws.Range('E30').Value = '  -2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0463'
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.21'
$ws.Range('D33').Value = '1.384.45'
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -6.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.51'
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.968'
$ws.Range('E36').Value = '  -4.59%  '
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.538'
$ws.Range('E39').Value = '  -3.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.817'
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.978'
$ws.Range('E42').Value = '  -2.66%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.18'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.76'
$ws.Range('E44').Value = '  -4.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.29'
$ws.Range('E45').Value = '  -4.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.20'
$ws.Range('E46').Value = '  -4.88%  '
$ws.Range('D47').Value = '1.718.77'
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.88'
$ws.Range('D49').Value = '0.0₆0100'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0970'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0498'
$ws.Range('E51').Value = '  -1.31%  '
